# Add 2022-Q1 data:
#  - The existing "总计" sheet becomes "2022-Q1" and is populated with the
#    new quarterly fund-holdings table.
#  - A brand-new "总计" sheet is appended at the end, re-using the old
#    totals table plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

$ws4 = $wb.Worksheets.Item("2021-Q4")
$wsQ1 = $wb.Worksheets.Item("总计")
$wsQ1.Name = "2022-Q1"

# ---------------------------------------------------------------------
# 1. Rebuild "2022-Q1" (the former "总计" sheet) with the fund table.
# ---------------------------------------------------------------------
$wsQ1.Cells.Clear()

# Reuse the header / index-column formatting already defined on 2021-Q4
# (style index carries bold font + border + centered alignment).
$ws4.Range("B1:H1").Copy()
$wsQ1.Range("B1:H1").PasteSpecial(-4122)
$ws4.Range("A2:A8").Copy()
$wsQ1.Range("A2:A12").PasteSpecial(-4122)

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $wsQ1.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$rows = @(
    @(0,  "700003", "平安策略先锋混合",                    "25.68", "72.64", "2.40", "0.6163", 8),
    @(1,  "004390", "平安转型创新灵活配置混合A",            "15.43", "88.10", "2.70", "0.4166", 9),
    @(2,  "012184", "大成创新趋势混合型证券投资基金A",      "8.05",  "72.67", "3.04", "0.2447", 10),
    @(3,  "004391", "平安转型创新灵活配置混合C",            "7.73",  "88.10", "2.70", "0.2087", 9),
    @(4,  "008274", "大成行业先锋混合A",                    "3.19",  "73.98", "3.04", "0.0970", 10),
    @(5,  "002945", "大成盛世精选灵活配置混合",             "1.52",  "70.48", "3.07", "0.0467", 10),
    @(6,  "011765", "兴银高端制造混合A",                    "1.01",  "93.23", "2.76", "0.0279", 6),
    @(7,  "008275", "大成行业先锋混合C",                    "0.45",  "73.98", "3.04", "0.0137", 10),
    @(8,  "011766", "兴银高端制造混合C",                    "0.39",  "93.23", "2.76", "0.0108", 6),
    @(9,  "012185", "大成创新趋势混合型证券投资基金C",      "0.14",  "72.67", "3.04", "0.0043", 10),
    @(10, "009027", "浦银安盛安远回报一年持有期混合A",      "0.79",  "20.03", "0.52", "0.0041", 10)
)

$r = 2
foreach ($row in $rows) {
    $wsQ1.Cells.Item($r, 1).Value = $row[0]

    # Text columns (B..G) must stay text, even though several look numeric
    # (e.g. fund codes with leading zeros, or "25.68"). Force text storage
    # via NumberFormat, then drop back to the Normal style so no stray
    # number-format style lingers on the cell.
    for ($col = 2; $col -le 7; $col++) {
        $cell = $wsQ1.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col - 1]
        $cell.Style = "Normal"
    }

    $wsQ1.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# ---------------------------------------------------------------------
# 2. Append the new "总计" sheet (old totals sheet, shifted + extra row).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Add($null, $wsQ1)
$wsTotal.Name = "总计"

$ws4.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$ws4.Range("A2:A3").Copy()
$wsTotal.Range("A2:A3").PasteSpecial(-4122)

$wsTotal.Cells.Item(1, 2).Value = "日期"
$wsTotal.Cells.Item(1, 3).Value = "持有数量(只)"
$wsTotal.Cells.Item(1, 4).Value = "持有市值(亿元)"

$wsTotal.Cells.Item(2, 1).Value = 0
$c = $wsTotal.Cells.Item(2, 2)
$c.NumberFormat = "@"
$c.Value = "2022-Q1"
$c.Style = "Normal"
$wsTotal.Cells.Item(2, 3).Value = 11
$wsTotal.Cells.Item(2, 4).Value = 1.69

$wsTotal.Cells.Item(3, 1).Value = 1
$c = $wsTotal.Cells.Item(3, 2)
$c.NumberFormat = "@"
$c.Value = "2021-Q4"
$c.Style = "Normal"
$wsTotal.Cells.Item(3, 3).Value = 7
$wsTotal.Cells.Item(3, 4).Value = 0.43

Write-Output "done"
